# Apply the "30_14_stimuli" worksheet update:
# - Fill in the new "carrier" column (D) for the practice + generic stimuli rows
# - Fill in the new "unique_video" / "unique_audio" pair-kind markers in column J
#   for the A/B practice-pair example rows
# - Fill in the new "unique_video" / "unique_audio" kind + carrier columns (C/D)
#   for the freshly-populated stimuli rows 14-21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("carrier") for rows 2-5 (practice p1-p4) and 6-9 (generic stimuli 1-4)
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"
$ws.Range("D6").Value = "can"
$ws.Range("D7").Value = "can"
$ws.Range("D8").Value = "do"
$ws.Range("D9").Value = "do"

# Column J ("pair_kind") for rows 6-9: unique video/audio markers
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# Rows 14-21: new "kind" (C) and "carrier" (D) values for stimuli 9-16
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
